# Add the three missing expense rows to the "data" sheet
# (row 3 was previously blank, rows 4-5 did not exist yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2890
$ws.Range("B3").Value = "чайник маме"

$ws.Range("A4").Value = 693
$ws.Range("B4").Value = "глобус"

$ws.Range("A5").Value = 500
$ws.Range("B5").Value = "подписка тинькоф и яндекс"
